$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- 1. Update the time_taken values (column F, rows 2-11) on the "data" sheet ---
$data.Cells.Item(2,6).Value = "2021-10-05 14:22:49.882442"
$data.Cells.Item(3,6).Value = "2021-10-05 14:22:49.882449"
$data.Cells.Item(4,6).Value = "2021-10-05 14:22:49.882453"
$data.Cells.Item(5,6).Value = "2021-10-05 14:22:49.882455"
$data.Cells.Item(6,6).Value = "2021-10-05 14:22:49.882458"
$data.Cells.Item(7,6).Value = "2021-10-05 14:22:49.882461"
$data.Cells.Item(8,6).Value = "2021-10-05 14:22:49.882463"
$data.Cells.Item(9,6).Value = "2021-10-05 14:22:49.882466"
$data.Cells.Item(10,6).Value = "2021-10-05 14:22:49.882468"
$data.Cells.Item(11,6).Value = "2021-10-05 14:22:49.882471"

# --- 2. Add the new "metadata" sheet, positioned after "data" ---
$meta = $wb.Worksheets.Add([Type]::Missing, $data)
$meta.Name = "metadata"

# Match the page margins used elsewhere in the workbook (values in points;
# 72pt = 1 inch).
$meta.PageSetup.LeftMargin = 54
$meta.PageSetup.RightMargin = 54
$meta.PageSetup.TopMargin = 72
$meta.PageSetup.BottomMargin = 72
$meta.PageSetup.HeaderMargin = 36
$meta.PageSetup.FooterMargin = 36

# Reuse the existing header style (bold font + border + center/top alignment)
# by copying the formatting from the "data" sheet header row, then overwrite
# the copied values with the metadata header labels.
$data.Range("B1:F1").Copy($meta.Range("B1:G1"))
$data.Range("A2").Copy($meta.Range("A2"))

$meta.Cells.Item(1,2).Value = "data_name"
$meta.Cells.Item(1,3).Value = "data_id"
$meta.Cells.Item(1,4).Value = "data_version"
$meta.Cells.Item(1,5).Value = "data_version_created"
$meta.Cells.Item(1,6).Value = "panel_query_time"
$meta.Cells.Item(1,7).Value = "panel_get_request"

# Row 2 - metadata values
$meta.Cells.Item(2,1).Value = 0
$meta.Cells.Item(2,2).Value = "Stickler syndrome"
$meta.Cells.Item(2,3).Value = 3
# data_version "2.22" must remain textual, not be coerced into a number
$meta.Cells.Item(2,4).NumberFormat = "@"
$meta.Cells.Item(2,4).Value = "2.22"
$meta.Cells.Item(2,5).Value = "2021-06-15T14:53:37.294988Z"
$meta.Cells.Item(2,6).Value = "2021-10-05 14:22:49.878748"
$meta.Cells.Item(2,7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/3/?format=json"

$meta.Range("A1").Select()
